$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.474
$ws.Range("A7").Value = -21.153
$ws.Range("B7").Value = 6.053
$ws.Range("B12").Value = 4.877000000000001
$ws.Range("E12").Value = 17.57
$ws.Range("D13").Value = -8.281000000000001
$ws.Range("D14").Value = -7.786999999999999
$ws.Range("B15").Value = 5.359000000000001
$ws.Range("A16").Value = -21.43
$ws.Range("D16").Value = -8.518000000000002
$ws.Range("D19").Value = -7.796000000000001
$ws.Range("A20").Value = -21.857
$ws.Range("B20").Value = 5.539
$ws.Range("B21").Value = 8.904
$ws.Range("B22").Value = 6.197
$ws.Range("D22").Value = -8.266999999999999
$ws.Range("E22").Value = 16.492
$ws.Range("B23").Value = 7.359999999999999
$ws.Range("A28").Value = -21.85
$ws.Range("A29").Value = -21.312
$ws.Range("B29").Value = 6.093
$ws.Range("E29").Value = 17.1
$ws.Range("A32").Value = -21.599
$ws.Range("B34").Value = 7.603999999999999
$ws.Range("E34").Value = 16.663
$ws.Range("D36").Value = -7.632
$ws.Range("A40").Value = -20.312
$ws.Range("B42").Value = 7.238000000000001
$ws.Range("B43").Value = 5.795
$ws.Range("E43").Value = 17.11
$ws.Range("B44").Value = 5.315
$ws.Range("B45").Value = 5.286
$ws.Range("A46").Value = -20.849
$ws.Range("B46").Value = 6.865
$ws.Range("D46").Value = -7.833
$ws.Range("E48").Value = 17.23
$ws.Range("B50").Value = 5.828
$ws.Range("D50").Value = -8.158000000000001
$ws.Range("A51").Value = -20.771
$ws.Range("B51").Value = 7.779000000000001
$ws.Range("A52").Value = -21.316
$ws.Range("A57").Value = -22.137
$ws.Range("A59").Value = -22.091
$ws.Range("E60").Value = 16.413
$ws.Range("A62").Value = -21.85
$ws.Range("A66").Value = -21.322
$ws.Range("B66").Value = 5.83
$ws.Range("B67").Value = 5.194999999999999
$ws.Range("E68").Value = 17.334
$ws.Range("E70").Value = 17.626
$ws.Range("A73").Value = -20.53
$ws.Range("E73").Value = 16.528
$ws.Range("A74").Value = -21.043
$ws.Range("B79").Value = 5.605
$ws.Range("B84").Value = 5.781000000000001
$ws.Range("E87").Value = 16.465
$ws.Range("A92").Value = -20.946
$ws.Range("B92").Value = 6.044
$ws.Range("E92").Value = 17.789
$ws.Range("D95").Value = -7.784000000000001
$ws.Range("B97").Value = 5.161
$ws.Range("D97").Value = -8.465
$ws.Range("A100").Value = -21.481
$ws.Range("E101").Value = 16.666
